$d = $word.ActiveDocument

# 1. Update the date in the letter header
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing-address paragraph ("3042 Pellier Pl, San Jose CA 95135")
#    into two lines ("3042 Pellier Pl" / "San Jose, CA 95135"), plus a new
#    trailing blank paragraph. Locate the paragraph by its distinctive text
#    so paragraph-index shifts elsewhere can't throw this off.
$addressPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*3042 Pellier Pl, San Jose CA 95135*") {
        $addressPara = $p
        break
    }
}
$addressPara.Range.Find.Execute(", San Jose CA 95135", $false, $false, $false, $false, $false,
                                 $true, 1, $false, "^pSan Jose, CA 95135^p", 2)

# 3. Remove the two empty paragraphs that directly follow the
#    "...Board of Directors" line (a No Spacing one, then the first of two
#    Title ones), leaving the remaining (second) Title paragraph in place.
#    Locate the paragraph by its distinctive text (not by index, since the
#    address-split above shifted everything after it by two).
$boardPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors*") {
        $boardPara = $p
        break
    }
}

$toDelete1 = $boardPara.Next()
$r1 = $d.Range($toDelete1.Range.Start, $toDelete1.Range.End)
$r1.Delete()

$toDelete2 = $boardPara.Next()
$r2 = $d.Range($toDelete2.Range.Start, $toDelete2.Range.End)
$r2.Delete()

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
